# [Kadastro App] Yeni kayit eklendi: 2946
# Appends one new record row to both the master "Kayitlar" sheet and the
# filtered "Erdemli" sheet (the new record's Birim is "Erdemli"), each as
# the next row right after their current last data row.

$wb = $excel.ActiveWorkbook

$recordNo   = "2946"
$tarih      = "2025-09-08"
$birim      = "Erdemli"
$parselSay  = "1"
$is         = "ÇAP"
$personel   = "CEMAL TİMUROĞLU (K.Teknisyeni)"

function Add-KayitRow($SheetName) {

    $ws = $wb.Worksheets.Item($SheetName)

    # Find the next empty row right after the current used range.
    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
    $newRow = $lastRow + 1

    # Use a leading apostrophe so numeric-looking values ("2946", the date,
    # "1") are stored as text, matching the existing columns in the sheet.
    $ws.Cells.Item($newRow, 1).Value = "'" + $recordNo
    $ws.Cells.Item($newRow, 2).Value = "'" + $tarih
    $ws.Cells.Item($newRow, 3).Value = $birim
    $ws.Cells.Item($newRow, 4).Value = "'" + $parselSay
    $ws.Cells.Item($newRow, 5).Value = $is
    $ws.Cells.Item($newRow, 6).Value = $personel
}

Add-KayitRow "Kayitlar"
Add-KayitRow "Erdemli"
